# Generate Report for Handoff
#
# The localization-status report is regenerated: the file
# "621bf619-44d3-44f4-9dcf-f140d569ffdb.md" moves from
# "Handed back: in sync with en-US" to "Ready for handoff", and its
# zh-cn / de-de "Latest Handoff Datetime" timestamps are refreshed to
# reflect the new handoff.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet: summary status for zh-cn (B) and de-de (C) columns ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- zh-cn sheet: per-language detail row for the same file ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $newStatus
$zhcn.Range("D3").Value = "2016-02-17 04:58:54"

# --- de-de sheet: per-language detail row for the same file ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $newStatus
$dede.Range("D3").Value = "2016-02-17 04:59:04"
